$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeRef, $text) {
    $ws.Range("ZZ1").NumberFormat = "@"
    $ws.Range("ZZ1").Value = $text
    $ws.Range("ZZ1").Copy()
    $ws.Range($rangeRef).PasteSpecial(-4163)
}

Set-TextValue "D2" "56.916.17"
Set-TextValue "E2" "  +1.16%  "
Set-TextValue "D3" "3.243.29"
Set-TextValue "E3" "  +0.08%  "
Set-TextValue "D5" "397.08"
Set-TextValue "E5" "  -1.12%  "
Set-TextValue "D6" "107.88"
Set-TextValue "E6" "  -2.98%  "
Set-TextValue "D7" "0.584"
Set-TextValue "E7" "  +4.59%  "
Set-TextValue "E9" "  -0.91%  "
Set-TextValue "D10" "39.35"
Set-TextValue "E10" "  -0.94%  "
Set-TextValue "D11" "0.0960"
Set-TextValue "E11" "  +6.80%  "
Set-TextValue "E12" "  +1.93%  "
Set-TextValue "D13" "3.754.88"
Set-TextValue "E13" "  -0.10%  "
Set-TextValue "D14" "8.31"
Set-TextValue "E14" "  +2.79%  "
Set-TextValue "E15" "  -1.74%  "
Set-TextValue "D16" "3.242.20"
Set-TextValue "E16" "  -0.49%  "
Set-TextValue "E17" "  -3.96%  "
Set-TextValue "E18" "  +3.85%  "
Set-TextValue "D19" "56.798.03"
Set-TextValue "E19" "  +1.12%  "
Set-TextValue "E20" "  -1.92%  "
Set-TextValue "D21" "0.0000111"
Set-TextValue "E21" "  +9.41%  "
Set-TextValue "E22" "  -1.82%  "
Set-TextValue "D23" "292.50"
Set-TextValue "E23" "  +1.03%  "
Set-TextValue "E24" "  -0.12%  "
Set-TextValue "E25" "  -2.53%  "
Set-TextValue "B26" "Filecoin"
Set-TextValue "C26" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D26" "8.02"
Set-TextValue "E26" "  -2.17%  "
Set-TextValue "B27" "EthereumClassic"
Set-TextValue "C27" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D27" "28.13"
Set-TextValue "E27" "  -0.41%  "
Set-TextValue "E28" "  -0.43%  "
Set-TextValue "E29" "  -0.80%  "
Set-TextValue "E30" "  -4.72%  "
Set-TextValue "E31" "  +0.04%  "
Set-TextValue "B32" "Hedera"
Set-TextValue "C32" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D32" "0.109"
Set-TextValue "E32" "  -2.35%  "
Set-TextValue "B33" "InjectiveProtocol"
Set-TextValue "C33" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D33" "41.19"
Set-TextValue "E33" "  +10.78%  "
Set-TextValue "E34" "  -1.76%  "
Set-TextValue "E35" "  -2.64%  "
Set-TextValue "E36" "  +0.94%  "
Set-TextValue "D37" "51.14"
Set-TextValue "E37" "  -0.44%  "
Set-TextValue "D38" "0.999"
Set-TextValue "E38" "  +0.02%  "
Set-TextValue "E39" "  -4.40%  "
Set-TextValue "E40" "  -3.56%  "
Set-TextValue "D41" "138.86"
Set-TextValue "E41" "  -0.46%  "
Set-TextValue "E42" "  +2.24%  "
Set-TextValue "D43" "3.95"
Set-TextValue "E43" "  -2.23%  "
Set-TextValue "E44" "  -3.15%  "
Set-TextValue "D45" "16.79"
Set-TextValue "E45" "  -1.46%  "
Set-TextValue "E46" "  -3.32%  "
Set-TextValue "E47" "  +8.32%  "
Set-TextValue "D48" "22.47"
Set-TextValue "E48" "  -1.15%  "
Set-TextValue "D49" "2.155.13"
Set-TextValue "E49" "  +0.39%  "
Set-TextValue "E50" "  -5.91%  "
Set-TextValue "E51" "  -7.27%  "

$ws.Range("ZZ1").Clear()
$excel.CutCopyMode = $false
